$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -1413928.571428571
$ws.Range("C3").Value = 5160000
$ws.Range("C4").Value = 2658928.571428571
$ws.Range("C5").Value = 8118214.285714286
$ws.Range("C6").Value = 4180000
$ws.Range("C7").Value = 18703214.28571429
